# RecordSet.sum_by handles case of summable key
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row 28: new "detail_a"/"detail_b" headers, right-aligned, no border
$ws.Range("L28").Value = "detail_a"
$ws.Range("M28").Value = "detail_b"
$ws.Range("O28").Value = "detail_a"
$ws.Range("P28").Value = "detail_b"
$ws.Range("L28:M28").HorizontalAlignment = -4152
$ws.Range("O28:P28").HorizontalAlignment = -4152

# Data rows 29-31: L/M detail_a/detail_b counters plus O/P summable-key sample
$ws.Range("L29").Value = 1
$ws.Range("M29").Value = 3
$ws.Range("O29").Value = 1
$ws.Range("P29").Value = 27

$ws.Range("L30").Value = 1
$ws.Range("M30").Value = 3
$ws.Range("O30").Value = 2
$ws.Range("P30").Value = 11

$ws.Range("L31").Value = 1
$ws.Range("M31").Value = 3
$ws.Range("O31").Value = 3
$ws.Range("P31").Value = 5

# Data rows 32-41: L/M detail_a/detail_b counters only
$ws.Range("L32").Value = 1
$ws.Range("M32").Value = 3

$ws.Range("L33").Value = 1
$ws.Range("M33").Value = 3

$ws.Range("L34").Value = 1
$ws.Range("M34").Value = 3

$ws.Range("L35").Value = 1
$ws.Range("M35").Value = 3

$ws.Range("L36").Value = 1
$ws.Range("M36").Value = 3

$ws.Range("L37").Value = 1
$ws.Range("M37").Value = 3

$ws.Range("L38").Value = 2
$ws.Range("M38").Value = 3

$ws.Range("L39").Value = 2
$ws.Range("M39").Value = 4

$ws.Range("L40").Value = 2
$ws.Range("M40").Value = 4

$ws.Range("L41").Value = 3
$ws.Range("M41").Value = 5

# Sort the O:P sample data range, which stamps a sortState on the sheet
$srt = $ws.Sort
$srt.SortFields.Clear()
$srt.SortFields.Add($ws.Range("N29:N45"))
$srt.SortFields.Add($ws.Range("O29:O45"))
$srt.SetRange($ws.Range("N29:O45"))
$srt.Header = 2
$srt.Apply()

# Restore the active selection to match the edited area
$ws.Range("O27").Select()
